# Adds a new "antenna_transport_cost" row to the SiteDevelopmentValues sheet,
# just above the existing "Complete/Full/Partial/Remote" infrastructure-scaling rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SiteDevelopmentValues")

# Copy the formatting of row 15 (A/B/D text style, C currency style) and insert it
# above row 16 -- this pushes the old rows 16-21 down to 17-22 while reusing the
# existing style records instead of minting new ones.
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(16).Insert()

# Fill in the new row's contents
$ws.Cells.Item(16, 1).Value = "antenna_transport_cost"
$ws.Cells.Item(16, 2).Value = "Cost to move antenna from mfg to site"
$ws.Cells.Item(16, 3).Value = 400000
$ws.Cells.Item(16, 4).Value = "from Bragg Crane quotes - about `$150k land voyage + `$250k sea voyage"

# This row is shorter than the copied row 15 (two lines of wrapped text instead of one)
$ws.Rows.Item(16).RowHeight = 36

# Update the selected cell to reflect the new active cell location
$ws.Range("C17").Select()
